$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "grouped_by_qtype": insert a new "Common-Divison" row at row 4 (the
# existing rows 4 "Multiplication" and 5 "Subtraction" shift down to 5 and 6),
# and refresh every metric value.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("grouped_by_qtype")

$ws1.Rows.Item(4).Insert()
$ws1.Rows.Item(4).ClearFormats()

$ws1.Range("B2").Value = 0.782051282051282
$ws1.Range("C2").Value = 0.8
$ws1.Range("D2").Value = 0.9721115537848606
$ws1.Range("E2").Value = 0.8776978417266188

$ws1.Range("B3").Value = 0.8356060606060606
$ws1.Range("C3").Value = 0.8452107279693487
$ws1.Range("D3").Value = 0.9865831842576028
$ws1.Range("E3").Value = 0.9104416013206768

$ws1.Range("A4").Value = "Common-Divison"
$ws1.Range("B4").Value = 0.875
$ws1.Range("C4").Value = 0.875
$ws1.Range("D4").Value = 1
$ws1.Range("E4").Value = 0.9333333333333333

$ws1.Range("A5").Value = "Multiplication"
$ws1.Range("B5").Value = 0.6666666666666666
$ws1.Range("C5").Value = 0.7182044887780549
$ws1.Range("D5").Value = 0.9028213166144201
$ws1.Range("E5").Value = 0.8

$ws1.Range("A6").Value = "Subtraction"
$ws1.Range("B6").Value = 0.7756591337099812
$ws1.Range("C6").Value = 0.7943587270973963
$ws1.Range("D6").Value = 0.9705449189985272
$ws1.Range("E6").Value = 0.873657695876972

# ---------------------------------------------------------------------------
# Sheet "grouped_by_distracted": values only, no structural change.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("grouped_by_distracted")

$ws2.Range("B2").Value = 0.7925
$ws2.Range("C2").Value = 0.8134462406979728
$ws2.Range("D2").Value = 0.9685304002444241
$ws2.Range("E2").Value = 0.8842398884239888

$ws2.Range("B3").Value = 0.75775
$ws2.Range("C3").Value = 0.7789771267026472
$ws2.Range("D3").Value = 0.9652866242038216
$ws2.Range("E3").Value = 0.8621817664628076

# ---------------------------------------------------------------------------
# Sheet "grouped_by_model": insert a new "deepseek" row at row 2 (existing
# "gemini", "llama3", "mistral" rows shift down by one), and refresh values.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("grouped_by_model")

$ws3.Rows.Item(2).Insert()
$ws3.Rows.Item(2).ClearFormats()

$ws3.Range("A2").Value = "deepseek"
$ws3.Range("B2").Value = 0.914
$ws3.Range("C2").Value = 0.944702842377261
$ws3.Range("D2").Value = 0.965662968832541
$ws3.Range("E2").Value = 0.955067920585162

$ws3.Range("A3").Value = "gemini"
$ws3.Range("B3").Value = 0.8555
$ws3.Range("C3").Value = 0.85678517776665
$ws3.Range("D3").Value = 0.9982497082847142
$ws3.Range("E3").Value = 0.9221234168687686

$ws3.Range("A4").Value = "llama3"
$ws3.Range("B4").Value = 0.4915
$ws3.Range("C4").Value = 0.5293484114162628
$ws3.Range("D4").Value = 0.8730017761989343
$ws3.Range("E4").Value = 0.659068052296346

$ws3.Range("A5").Value = "mistral"
$ws3.Range("B5").Value = 0.8395
$ws3.Range("C5").Value = 0.83991995997999
$ws3.Range("D5").Value = 0.9994047619047619
$ws3.Range("E5").Value = 0.9127480293558032

# ---------------------------------------------------------------------------
# Sheet "confusion_matrix": values only, no structural change.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("confusion_matrix")

$ws4.Range("B2").Value = 103
$ws4.Range("C2").Value = 727
$ws4.Range("D2").Value = 3170

$ws4.Range("B3").Value = 109
$ws4.Range("C3").Value = 860
$ws4.Range("D3").Value = 3031
